$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-17 11:07:41"
$wsZh.Range("H3").Value = "2016-03-17 11:07:58"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-17 11:07:44"
$wsDe.Range("H3").Value = "2016-03-17 11:08:05"
